# Saldo.xlsx update — applies the row-level changes described by the commit:
#   - row 6  (MARCUS / 000834301): balance 102114.63 -> 104229.07
#   - row 7  (DOUGLAS / 004384167 / 100000): row removed entirely
#   - two new rows inserted after the (now shifted) MARCUS/004752519 row:
#       004452476 IVONE  40626.84
#       004399832 EULER  38911.4
#   - one new row inserted before SURAMA (004205505):
#       003894173 ANDREA 1309.01
#   - old row for IVONE / 004452476 / 713.31 (between ELEUSE and FABIANO) removed
#
# Operations are applied from the bottom of the sheet upward so that earlier
# (lower) row numbers used below are never invalidated by a later edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Remove the old IVONE (004452476) row that sits between ELEUSE (717.24)
#        and FABIANO (680.57), currently row 42.
$ws.Rows(42).Delete()

# --- 2) Insert a new row for ANDREA (003894173 / 1309.01) right before SURAMA
#        (004205505), currently row 13.
$ws.Rows(13).Insert()
$ws.Range("A13").Value = "'003894173"
$ws.Range("B13").Value = "ANDREA"
$ws.Range("C13").Value = 1309.01

# --- 3) Insert two new rows right before THOMAS (004224011), currently row 9:
#        final order must read IVONE, EULER, THOMAS - since each Insert() at
#        the same row number pushes the previously-written row back down, the
#        row written LAST ends up on top, so EULER is written first, IVONE
#        second.
$ws.Rows(9).Insert()
$ws.Range("A9").Value = "'004399832"
$ws.Range("B9").Value = "EULER"
$ws.Range("C9").Value = 38911.4

$ws.Rows(9).Insert()
$ws.Range("A9").Value = "'004452476"
$ws.Range("B9").Value = "IVONE"
$ws.Range("C9").Value = 40626.84

# --- 4) Remove the DOUGLAS (004384167 / 100000) row, currently row 7.
$ws.Rows(7).Delete()

# --- 5) Update MARCUS (000834301) balance on row 6.
$ws.Range("C6").Value = 104229.07
